# Update the form title on the "settings" sheet, then make that sheet the
# active tab (with a fresh selection), matching the recorded edit.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "Invest, Insure?"

$settings.Activate()
$settings.Range("A2").Select()
